$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    Target shape:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Discover the action-packed gameplay of Battleship
#                    Direct Hit Megaways. Read our review now and play
#                    for free.</w:t></w:r>
#      </w:p>
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"          # drop inherited Heading1 pStyle

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover the action-packed gameplay of Battleship Direct Hit Megaways. Read our review now and play for free.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) At the end of the document: drop the bold "Play Battleship Direct Hit
#    for Free - Read Our Review Now" paragraph entirely, and rewrite the
#    italic paragraph's text (keeping its <w:i/> formatting).
# ---------------------------------------------------------------------------

$oldText = "Discover the action-packed gameplay of Battleship Direct Hit Megaways. Read our review now and play for free."
$titleText = "Play Battleship Direct Hit for Free - Read Our Review Now"

# Locate the trailing duplicate-title (bold) paragraph and the following
# italic paragraph by content, rather than assuming fixed indices.
$boldPara = $null
$italicPara = $null
$cnt = $d.Paragraphs.Count
for ($i = $cnt; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($italicPara -eq $null) -and ($t -like "*$oldText*")) {
        $italicPara = $p
    } elseif (($italicPara -ne $null) -and ($boldPara -eq $null) -and ($t -like "*$titleText*")) {
        $boldPara = $p
        break
    }
}

$killRange = $d.Range($boldPara.Range.Start, $italicPara.Range.Start)
$killRange.Delete()

$cnt = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($cnt)
$newText = 'Please create an image featuring a happy Maya warrior wearing glasses for the game "Battleship Direct Hit". The image should be in a cartoon style and should capture the essence of the game''s naval battle theme in a fun and engaging way. It should also feature the game''s title prominently. Be creative and use vibrant colors and dynamic imagery to attract players to this exciting slot game.'

$start = $italicPara.Range.Start
$textRange = $d.Range($start, $start + $oldText.Length)
$textRange.Text = $newText

Write-Host "Done."
